$d = $word.ActiveDocument

$d.Content.Find.Execute("606÷9=67, 3", $true, $false, $false, $false, $false, $true, 1, $false, "308÷5=61, 3", 2) | Out-Null
$d.Content.Find.Execute("657÷8=82, 1", $true, $false, $false, $false, $false, $true, 1, $false, "820÷6=136, 4", 2) | Out-Null
$d.Content.Find.Execute("735÷3=245, 0", $true, $false, $false, $false, $false, $true, 1, $false, "731÷4=182, 3", 2) | Out-Null
$d.Content.Find.Execute("639÷7=91, 2", $true, $false, $false, $false, $false, $true, 1, $false, "536÷4=134, 0", 2) | Out-Null
$d.Content.Find.Execute("532÷9=59, 1", $true, $false, $false, $false, $false, $true, 1, $false, "734÷6=122, 2", 2) | Out-Null
$d.Content.Find.Execute("201÷5=40, 1", $true, $false, $false, $false, $false, $true, 1, $false, "794÷9=88, 2", 2) | Out-Null
$d.Content.Find.Execute("916÷2=458, 0", $true, $false, $false, $false, $false, $true, 1, $false, "550÷8=68, 6", 2) | Out-Null
$d.Content.Find.Execute("322÷8=40, 2", $true, $false, $false, $false, $false, $true, 1, $false, "382÷3=127, 1", 2) | Out-Null
$d.Content.Find.Execute("602÷8=75, 2", $true, $false, $false, $false, $false, $true, 1, $false, "463÷4=115, 3", 2) | Out-Null
$d.Content.Find.Execute("209÷5=41, 4", $true, $false, $false, $false, $false, $true, 1, $false, "370÷9=41, 1", 2) | Out-Null
$d.Content.Find.Execute("649÷5=129, 4", $true, $false, $false, $false, $false, $true, 1, $false, "711÷5=142, 1", 2) | Out-Null
$d.Content.Find.Execute("503÷5=100, 3", $true, $false, $false, $false, $false, $true, 1, $false, "295÷7=42, 1", 2) | Out-Null
$d.Content.Find.Execute("671÷9=74, 5", $true, $false, $false, $false, $false, $true, 1, $false, "378÷4=94, 2", 2) | Out-Null
$d.Content.Find.Execute("554÷2=277, 0", $true, $false, $false, $false, $false, $true, 1, $false, "613÷8=76, 5", 2) | Out-Null
$d.Content.Find.Execute("140÷9=15, 5", $true, $false, $false, $false, $false, $true, 1, $false, "126÷6=21, 0", 2) | Out-Null
$d.Content.Find.Execute("371÷2=185, 1", $true, $false, $false, $false, $false, $true, 1, $false, "758÷7=108, 2", 2) | Out-Null
$d.Content.Find.Execute("577÷4=144, 1", $true, $false, $false, $false, $false, $true, 1, $false, "764÷2=382, 0", 2) | Out-Null
$d.Content.Find.Execute("748÷6=124, 4", $true, $false, $false, $false, $false, $true, 1, $false, "118÷4=29, 2", 2) | Out-Null
$d.Content.Find.Execute("699÷6=116, 3", $true, $false, $false, $false, $false, $true, 1, $false, "645÷6=107, 3", 2) | Out-Null
$d.Content.Find.Execute("686÷8=85, 6", $true, $false, $false, $false, $false, $true, 1, $false, "789÷6=131, 3", 2) | Out-Null
$d.Content.Find.Execute("829÷7=118, 3", $true, $false, $false, $false, $false, $true, 1, $false, "465÷8=58, 1", 2) | Out-Null
$d.Content.Find.Execute("922÷9=102, 4", $true, $false, $false, $false, $false, $true, 1, $false, "636÷7=90, 6", 2) | Out-Null
$d.Content.Find.Execute("474÷6=79, 0", $true, $false, $false, $false, $false, $true, 1, $false, "491÷7=70, 1", 2) | Out-Null
$d.Content.Find.Execute("594÷6=99, 0", $true, $false, $false, $false, $false, $true, 1, $false, "826÷3=275, 1", 2) | Out-Null
$d.Content.Find.Execute("266÷5=53, 1", $true, $false, $false, $false, $false, $true, 1, $false, "568÷9=63, 1", 2) | Out-Null
